$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 205.33333
$ws.Range("I11").Value = 205.33333
$ws.Range("K11").Value = 205.33333
$ws.Range("M11").Value = -65.33332999999999

$ws.Range("H76").Value = 1565390.5
$ws.Range("I76").Value = 2345172.2
$ws.Range("J76").Value = 5827
$ws.Range("K76").Value = 2345172.2
$ws.Range("L76").Value = 5827
$ws.Range("M76").Value = -2344857.2
$ws.Range("N76").Value = -6457

$ws.Range("H79").Value = 1565390.5
$ws.Range("I79").Value = 2345172.2
$ws.Range("J79").Value = 5827
$ws.Range("K79").Value = 2345172.2
$ws.Range("L79").Value = 5827
$ws.Range("M79").Value = -2344080.2
$ws.Range("N79").Value = -8011

$ws.Range("H113").Value = 24700
$ws.Range("I113").Value = 24700
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 24700
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -21446
$ws.Range("N113").ClearContents() | Out-Null

$ws.Range("H129").Value = 903.9846
$ws.Range("I129").Value = 1219.8
$ws.Range("K129").Value = 3659.4
$ws.Range("M129").Value = 1340.6

$ws.Range("H133").Value = 79999
$ws.Range("J133").Value = 79999
$ws.Range("L133").Value = 79999
$ws.Range("N133").Value = -90119

$ws.Range("H137").Value = 1653.2
$ws.Range("I137").Value = 1345.2727
$ws.Range("K137").Value = 4035.8181
$ws.Range("M137").Value = -1485.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2694.9285
$ws.Range("I61").Value = 2192.2563
$ws.Range("J61").Value = 9229.666999999999
$ws.Range("K61").Value = 2192.2563
$ws.Range("L61").Value = 9229.666999999999
$ws.Range("M61").Value = -1980.2563
$ws.Range("N61").Value = -9653.666999999999

$ws.Range("H74").Value = 1232.2858
$ws.Range("I74").Value = 981.17645
$ws.Range("J74").Value = 1620.3636
$ws.Range("K74").Value = 981.17645
$ws.Range("L74").Value = 1620.3636
$ws.Range("M74").Value = -107.17645
$ws.Range("N74").Value = -3368.3636

$ws.Range("H77").Value = 1232.2858
$ws.Range("I77").Value = 981.17645
$ws.Range("J77").Value = 1620.3636
$ws.Range("K77").Value = 4905.882250000001
$ws.Range("L77").Value = 8101.817999999999
$ws.Range("M77").Value = -537.8822500000006
$ws.Range("N77").Value = -16837.818

$ws.Range("H88").Value = 23769.6
$ws.Range("J88").Value = 28974.625
$ws.Range("L88").Value = 28974.625
$ws.Range("N88").Value = -29786.625

$ws.Range("H91").Value = 23769.6
$ws.Range("J91").Value = 28974.625
$ws.Range("L91").Value = 28974.625
$ws.Range("N91").Value = -31782.625

$ws.Range("H110").Value = 1772.7
$ws.Range("I110").Value = 1428.0769
$ws.Range("K110").Value = 1428.0769
$ws.Range("M110").Value = 616.9231

$ws.Range("H122").Value = 4049.875
$ws.Range("I122").Value = 4879.8
$ws.Range("K122").Value = 14639.4
$ws.Range("M122").Value = -12189.4

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents() | Out-Null

$ws.Range("H136").Value = 2694.9285
$ws.Range("I136").Value = 2192.2563
$ws.Range("J136").Value = 9229.666999999999
$ws.Range("K136").Value = 6576.7689
$ws.Range("L136").Value = 27689.001
$ws.Range("M136").Value = -4026.7689
$ws.Range("N136").Value = -32789.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 86135.66
$ws.Range("I86").Value = 1400.8235
$ws.Range("J86").Value = 246190.33
$ws.Range("K86").Value = 1400.8235
$ws.Range("L86").Value = 246190.33
$ws.Range("M86").Value = -277.8235
$ws.Range("N86").Value = -248436.33

$ws.Range("H89").Value = 86135.66
$ws.Range("I89").Value = 1400.8235
$ws.Range("J89").Value = 246190.33
$ws.Range("K89").Value = 7004.1175
$ws.Range("L89").Value = 1230951.65
$ws.Range("M89").Value = -1388.1175
$ws.Range("N89").Value = -1242183.65

$ws.Range("H134").Value = 6199.1914
$ws.Range("I134").Value = 5970.683
$ws.Range("J134").Value = 7760.6665
$ws.Range("K134").Value = 17912.049
$ws.Range("L134").Value = 23281.9995
$ws.Range("M134").Value = -15377.049
$ws.Range("N134").Value = -28351.9995

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 739.6667
$ws.Range("I16").Value = 707.125
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 707.125
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -420.125
$ws.Range("N16").Value = -1574

$ws.Range("H23").Value = 63138.332
$ws.Range("I23").Value = 49900
$ws.Range("K23").Value = 49900
$ws.Range("M23").Value = -49660

$ws.Range("H27").Value = 63138.332
$ws.Range("I27").Value = 49900
$ws.Range("K27").Value = 49900
$ws.Range("M27").Value = -49708

$ws.Range("H31").Value = 2544.8044
$ws.Range("I31").Value = 1641.6451
$ws.Range("J31").Value = 4411.3335
$ws.Range("K31").Value = 1641.6451
$ws.Range("L31").Value = 4411.3335
$ws.Range("M31").Value = -1346.6451
$ws.Range("N31").Value = -5001.3335

$ws.Range("H34").Value = 2544.8044
$ws.Range("I34").Value = 1641.6451
$ws.Range("J34").Value = 4411.3335
$ws.Range("K34").Value = 1641.6451
$ws.Range("L34").Value = 4411.3335
$ws.Range("M34").Value = -1439.6451
$ws.Range("N34").Value = -4815.3335

$ws.Range("H105").Value = 1323.5
$ws.Range("I105").Value = 1225.2858
$ws.Range("K105").Value = 1225.2858
$ws.Range("M105").Value = 521.7141999999999

$ws.Range("H113").Value = 739.6667
$ws.Range("I113").Value = 707.125
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 707.125
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1462.875
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 213747.69
$ws.Range("I4").Value = 21845.217
$ws.Range("J4").Value = 1685000
$ws.Range("K4").Value = 65535.651
$ws.Range("L4").Value = 5055000
$ws.Range("M4").Value = -65423.651
$ws.Range("N4").Value = -5055224

$ws.Range("H5").Value = 707.8
$ws.Range("J5").Value = 745.6667
$ws.Range("L5").Value = 2237.0001
$ws.Range("N5").Value = -2461.0001

$ws.Range("H11").Value = 1094.25
$ws.Range("I11").Value = 700
$ws.Range("J11").Value = 1488.5
$ws.Range("K11").Value = 2100
$ws.Range("L11").Value = 4465.5
$ws.Range("M11").Value = -1960
$ws.Range("N11").Value = -4745.5

$ws.Range("H107").Value = 880.8125
$ws.Range("J107").Value = 880.8125
$ws.Range("L107").Value = 2642.4375
$ws.Range("N107").Value = -6482.4375

$ws.Range("H121").Value = 704.6
$ws.Range("J121").Value = 842.8570999999999
$ws.Range("L121").Value = 2528.5713
$ws.Range("N121").Value = -5148.5713

$ws.Range("H122").Value = 961.7273
$ws.Range("J122").Value = 1022.5
$ws.Range("L122").Value = 9202.5
$ws.Range("N122").Value = -14102.5

$ws.Range("H131").Value = 10626.226
$ws.Range("J131").Value = 11895.841
$ws.Range("L131").Value = 35687.523
$ws.Range("N131").Value = -45767.523

$ws.Range("H135").Value = 707.8
$ws.Range("J135").Value = 745.6667
$ws.Range("L135").Value = 6711.0003
$ws.Range("N135").Value = -11781.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 65670.336
$ws.Range("I19").Value = 56999
$ws.Range("J19").Value = 70006
$ws.Range("K19").Value = 56999
$ws.Range("L19").Value = 70006
$ws.Range("M19").Value = -56711
$ws.Range("N19").Value = -70582

$ws.Range("H132").Value = 1070751.8
$ws.Range("I132").Value = 1327976.8
$ws.Range("J132").Value = 5105.7144
$ws.Range("K132").Value = 3983930.4
$ws.Range("L132").Value = 15317.1432
$ws.Range("M132").Value = -3981400.4
$ws.Range("N132").Value = -20377.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4833.1665
$ws.Range("I61").Value = 5249.75
$ws.Range("K61").Value = 5249.75
$ws.Range("M61").Value = -5047.75

$ws.Range("H113").Value = 4833.1665
$ws.Range("I113").Value = 5249.75
$ws.Range("K113").Value = 5249.75
$ws.Range("M113").Value = -3079.75

$ws.Range("H122").Value = 6750.625
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws.Range("H128").Value = 49999.5
$ws.Range("J128").Value = 49999.5
$ws.Range("L128").Value = 49999.5
$ws.Range("N128").Value = -59959.5

$ws.Range("H132").Value = 1827.0938
$ws.Range("I132").Value = 1611.1177
$ws.Range("K132").Value = 4833.3531
$ws.Range("M132").Value = -2303.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 236097
$ws.Range("I122").Value = 313962.66
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 941887.98
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -939437.98
$ws.Range("N122").Value = -12400

$ws.Range("H132").Value = 2148.5
$ws.Range("I132").Value = 1768.5883
$ws.Range("J132").Value = 3071.1428
$ws.Range("K132").Value = 5305.7649
$ws.Range("L132").Value = 9213.428400000001
$ws.Range("M132").Value = -2775.7649
$ws.Range("N132").Value = -14273.4284

$ws.Range("H136").Value = 14621836
$ws.Range("I136").Value = 19842958
$ws.Range("K136").Value = 59528874
$ws.Range("M136").Value = -59526324

$ws.Range("H140").Value = 57000
$ws.Range("J140").Value = 57000
$ws.Range("L140").Value = 57000
$ws.Range("N140").Value = -67360

$ws.Range("H141").Value = 73749.914
$ws.Range("J141").Value = 73749.914
$ws.Range("L141").Value = 73749.914
$ws.Range("N141").Value = -84109.914
